$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume change) per upstream data refresh.

$ws.Range("D2").Value = '42.902.99'
$ws.Range("E2").Value = '  -1.45%  '

$ws.Range("D3").Value = '2.340.37'
$ws.Range("E3").Value = '  +1.06%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.87'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.51%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.63'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.98%  '

$ws.Range("E7").Value = '  -5.13%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.511'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -4.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.98'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.59%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.16'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.15%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0800'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.11%  '

$ws.Range("E13").Value = '  -0.50%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.83'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.87'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +5.58%  '

$ws.Range("D16").Value = '2.339.89'
$ws.Range("E16").Value = '  +1.31%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.803'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.13%  '

$ws.Range("D18").Value = '42.820.39'
$ws.Range("E18").Value = '  -1.42%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.23'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.15%  '

$ws.Range("D20").Value = '0.0₃0911'
$ws.Range("E20").Value = '  -2.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.71'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -6.21%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.93'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.54%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.02'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.12%  '

$ws.Range("E24").Value = '  -2.36%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.56'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.37%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.60'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +3.27%  '

$ws.Range("E28").Value = '  +9.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.98'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -5.70%  '

$ws.Range("E30").Value = '  -3.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '159.96'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -4.93%  '

$ws.Range("E32").Value = '  -0.10%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.13'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.66'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +7.51%  '

$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.45'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -3.02%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0727'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.58%  '

$ws.Range("B37").Value = 'Celestia'
$ws.Range("C37").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.32'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.96'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -5.12%  '

$ws.Range("E39").Value = '  -0.64%  '

$ws.Range("E40").Value = '  -3.11%  '

$ws.Range("E41").Value = '  -2.96%  '

$ws.Range("E42").Value = '  +1.76%  '

$ws.Range("D43").Value = '2.027.39'
$ws.Range("E43").Value = '  +2.82%  '

$ws.Range("E44").Value = '  -1.66%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.85'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.89%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.34'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +4.09%  '

$ws.Range("E47").Value = '  -2.13%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '56.28'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.12%  '

$ws.Range("E49").Value = '  -1.18%  '

$ws.Range("D50").Value = '2.565.97'
$ws.Range("E50").Value = '  +0.92%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.66'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.26%  '
